$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.068.35'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.28%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.832.28'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.09%  '

# Row 4
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.30%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '243.61'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.96%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6287'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.62%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.003'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.22%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07469'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -1.79%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2931'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.92%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.01'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +1.56%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07713'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.51%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.825.34'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.30%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.006'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.14%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6664'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.76%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '83.23'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.98%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.000009427'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -2.86%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.047'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.14%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '29.068.11'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.29%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.59'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +2.25%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '224.28'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.31%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.003'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.30%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.106'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.15%  '

# Row 23
$ws.Range('B23').NumberFormat = '@'
$ws.Range('B23').Value = 'BinanceUSD'
$ws.Range('C23').NumberFormat = '@'
$ws.Range('C23').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.003'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.29%  '

# Row 24
$ws.Range('B24').NumberFormat = '@'
$ws.Range('B24').Value = 'Monero'
$ws.Range('C24').NumberFormat = '@'
$ws.Range('C24').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '159.94'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.11%  '

# Row 25
$ws.Range('B25').NumberFormat = '@'
$ws.Range('B25').Value = 'Stellar'
$ws.Range('C25').NumberFormat = '@'
$ws.Range('C25').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1404'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +2.74%  '

# Row 26
$ws.Range('B26').NumberFormat = '@'
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').NumberFormat = '@'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.501'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +1.37%  '

# Row 27
$ws.Range('B27').NumberFormat = '@'
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').NumberFormat = '@'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.93'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.53%  '

# Row 28
$ws.Range('B28').NumberFormat = '@'
$ws.Range('B28').Value = 'PancakeSwap'
$ws.Range('C28').NumberFormat = '@'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.500'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.87%  '

# Row 29
$ws.Range('B29').NumberFormat = '@'
$ws.Range('B29').Value = 'Filecoin'
$ws.Range('C29').NumberFormat = '@'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.126'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +1.84%  '

# Row 30
$ws.Range('B30').NumberFormat = '@'
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').NumberFormat = '@'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.057'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.90%  '

# Row 31
$ws.Range('B31').NumberFormat = '@'
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').NumberFormat = '@'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.05468'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +5.57%  '

# Row 32
$ws.Range('B32').NumberFormat = '@'
$ws.Range('B32').Value = 'Toncoin'
$ws.Range('C32').NumberFormat = '@'
$ws.Range('C32').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.198'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.59%  '

# Row 33
$ws.Range('B33').NumberFormat = '@'
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').NumberFormat = '@'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7504'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +1.64%  '

# Row 34
$ws.Range('B34').NumberFormat = '@'
$ws.Range('B34').Value = 'LidoDAOToken'
$ws.Range('C34').NumberFormat = '@'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.853'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.75%  '

# Row 35
$ws.Range('B35').NumberFormat = '@'
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').NumberFormat = '@'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.137'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.48%  '

# Row 36
$ws.Range('B36').NumberFormat = '@'
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').NumberFormat = '@'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.615'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -3.00%  '

# Row 37
$ws.Range('B37').NumberFormat = '@'
$ws.Range('B37').Value = 'Maker'
$ws.Range('C37').NumberFormat = '@'
$ws.Range('C37').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.231.42'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -2.48%  '

# Row 38
$ws.Range('B38').NumberFormat = '@'
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').NumberFormat = '@'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.752'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.27%  '

# Row 39
$ws.Range('B39').NumberFormat = '@'
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').NumberFormat = '@'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01786'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.19%  '

# Row 40
$ws.Range('B40').NumberFormat = '@'
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').NumberFormat = '@'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.642'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +7.00%  '

# Row 41
$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8956'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.28%  '

# Row 42
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.003'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.17%  '

# Row 43
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '101.76'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.28%  '

# Row 44
$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'BabyDogeCoin'
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.00000000124'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +1.31%  '

# Row 45
$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '65.63'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +2.29%  '

# Row 46
$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value = 'Mantle'
$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5097'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.21%  '

# Row 47
$ws.Range('B47').NumberFormat = '@'
$ws.Range('B47').Value = 'TheSandbox'
$ws.Range('C47').NumberFormat = '@'
$ws.Range('C47').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4041'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +1.66%  '

# Row 48
$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.929'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +1.20%  '

# Row 49
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'XinFinNetwork'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.07262'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +6.59%  '

# Row 50
$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.658'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +2.94%  '

# Row 51
$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05803'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.94%  '
